$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Unique"
$ws.Range("B7").Value = "Akash Shahapure"
$ws.Range("C7").Value = "akash.shahapure@haqdarshak.com"
$ws.Range("D7").Value = 4.836565017700195
$ws.Range("E7").Value = "03/07/2024 16:55:38"
$ws.Range("F7").Value = "03/07/2024 16:58:56"
$ws.Range("G7").Value = 198
$ws.Range("H7").Value = "cases_report_Tata Power (Maithon Power Ltd)_(All States)_2024-06-19"
$ws.Range("I7").Value = 17916
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 393

$ws.Range("A8").Value = "All"
$ws.Range("B8").Value = "Akash Shahapure"
$ws.Range("C8").Value = "akash.shahapure@haqdarshak.com"
$ws.Range("D8").Value = 4.836565017700195
$ws.Range("E8").Value = "03/07/2024 16:55:38"
$ws.Range("F8").Value = "03/07/2024 16:59:06"
$ws.Range("G8").Value = 208
$ws.Range("H8").Value = "cases_report_Tata Power (Maithon Power Ltd)_(All States)_2024-06-19"
$ws.Range("I8").Value = 17917
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 393
